$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("S6 Table")

# Update the B, C, D raw data values for rows 4-32 (E, F, G are formulas that recalc automatically)
$ws.Cells.Item(4, 2).Value = 0.48030739673391
$ws.Cells.Item(4, 3).Value = 0.767326732673267
$ws.Cells.Item(4, 4).Value = 1.80459770114943
$ws.Cells.Item(5, 2).Value = 0.468025949953661
$ws.Cells.Item(5, 3).Value = 0.742496050552923
$ws.Cells.Item(5, 4).Value = 1.90551181102362
$ws.Cells.Item(6, 2).Value = 4.3
$ws.Cells.Item(6, 3).Value = 6.84615384615385
$ws.Cells.Item(6, 4).Value = 11.1333333333333
$ws.Cells.Item(7, 2).Value = 1.16239316239316
$ws.Cells.Item(7, 3).Value = 1.24303797468354
$ws.Cells.Item(7, 4).Value = 3.90350877192983
$ws.Cells.Item(8, 2).Value = 0.625954198473283
$ws.Cells.Item(8, 3).Value = 0.680426098535286
$ws.Cells.Item(8, 4).Value = 2.24651162790698
$ws.Cells.Item(9, 2).Value = 4.98666666666667
$ws.Cells.Item(9, 3).Value = 6.73529411764706
$ws.Cells.Item(9, 4).Value = 13.1612903225806
$ws.Cells.Item(10, 2).Value = 0.89048473967684
$ws.Cells.Item(10, 3).Value = 1.01934235976789
$ws.Cells.Item(10, 4).Value = 3.1864406779661
$ws.Cells.Item(11, 2).Value = 0.361987911349899
$ws.Cells.Item(11, 3).Value = 0.422578184591915
$ws.Cells.Item(11, 4).Value = 1.62595419847328
$ws.Cells.Item(12, 2).Value = 0.506711409395973
$ws.Cells.Item(12, 3).Value = 0.619426751592357
$ws.Cells.Item(12, 4).Value = 2.8062015503876
$ws.Cells.Item(13, 2).Value = 0.37424789410349
$ws.Cells.Item(13, 3).Value = 0.473365617433414
$ws.Cells.Item(13, 4).Value = 2.8062015503876
$ws.Cells.Item(14, 2).Value = 0.594360086767896
$ws.Cells.Item(14, 3).Value = 0.674460431654676
$ws.Cells.Item(14, 4).Value = 2.8062015503876
$ws.Cells.Item(15, 2).Value = 0.598654708520179
$ws.Cells.Item(15, 3).Value = 0.678832116788321
$ws.Cells.Item(15, 4).Value = 2.8062015503876
$ws.Cells.Item(16, 2).Value = 0.505862646566164
$ws.Cells.Item(16, 3).Value = 0.619426751592357
$ws.Cells.Item(16, 4).Value = 2.8062015503876
$ws.Cells.Item(17, 2).Value = 0.485345838218054
$ws.Cells.Item(17, 3).Value = 0.525686977299881
$ws.Cells.Item(17, 4).Value = 2.064
$ws.Cells.Item(18, 2).Value = 0.521345407503234
$ws.Cells.Item(18, 3).Value = 0.537794299876084
$ws.Cells.Item(18, 4).Value = 2.064
$ws.Cells.Item(19, 2).Value = 0.584229390681004
$ws.Cells.Item(19, 3).Value = 0.735751295336788
$ws.Cells.Item(19, 4).Value = 2.91558441558442
$ws.Cells.Item(20, 2).Value = 0.585106382978723
$ws.Cells.Item(20, 3).Value = 0.737478411053541
$ws.Cells.Item(20, 4).Value = 2.84177215189873
$ws.Cells.Item(21, 2).Value = 0.631274131274131
$ws.Cells.Item(21, 3).Value = 0.784277879341865
$ws.Cells.Item(21, 4).Value = 2.775
$ws.Cells.Item(22, 2).Value = 0.576298701298701
$ws.Cells.Item(22, 3).Value = 0.697784810126582
$ws.Cells.Item(22, 4).Value = 2.13114754098361
$ws.Cells.Item(23, 2).Value = 0.578352180936995
$ws.Cells.Item(23, 3).Value = 0.710191082802548
$ws.Cells.Item(23, 4).Value = 2.14814814814815
$ws.Cells.Item(24, 2).Value = 0.615780445969125
$ws.Cells.Item(24, 3).Value = 0.732348111658456
$ws.Cells.Item(24, 4).Value = 2.07539682539683
$ws.Cells.Item(25, 2).Value = 0.5
$ws.Cells.Item(25, 3).Value = 0.57487922705314
$ws.Cells.Item(25, 4).Value = 1.25917431192661
$ws.Cells.Item(26, 2).Value = 0.483122362869198
$ws.Cells.Item(26, 3).Value = 0.551671732522796
$ws.Cells.Item(26, 4).Value = 1.18105263157895
$ws.Cells.Item(27, 2).Value = 0.506912442396313
$ws.Cells.Item(27, 3).Value = 0.590984974958264
$ws.Cells.Item(27, 4).Value = 1.30788177339902
$ws.Cells.Item(28, 2).Value = 0.671186440677966
$ws.Cells.Item(28, 3).Value = 0.804878048780488
$ws.Cells.Item(28, 4).Value = 1.53281853281853
$ws.Cells.Item(29, 2).Value = 0.497797356828194
$ws.Cells.Item(29, 3).Value = 0.576923076923077
$ws.Cells.Item(29, 4).Value = 1.25507900677201
$ws.Cells.Item(30, 2).Value = 0.949333333333333
$ws.Cells.Item(30, 3).Value = 1.15263157894737
$ws.Cells.Item(30, 4).Value = 4.67924528301887
$ws.Cells.Item(31, 2).Value = 0.484057971014493
$ws.Cells.Item(31, 3).Value = 0.512135922330097
$ws.Cells.Item(31, 4).Value = 1.56569343065693
$ws.Cells.Item(32, 2).Value = 0.470170454545455
$ws.Cells.Item(32, 3).Value = 0.498817966903073
$ws.Cells.Item(32, 4).Value = 1.62406015037594

# Recalculate formulas (E, F, G columns and MAX/MIN rows) to reflect updated data
$excel.Calculate()

# Update the sheet selection to match the committed state (single cell B4 selected instead of B4:B32)
$ws.Range("B4").Select()
